$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "22.429.76"
$ws.Range("E2").Value = "  +0.07%  "

# Row 3
$ws.Range("D3").Value = "1.571.07"
$ws.Range("E3").Value = "  +0.58%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("E5").Value = "  +0.05%  "

# Row 6
$ws.Range("D6").Value = "'287.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.25%  "

# Row 7
$ws.Range("D7").Value = "'0.3688"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.16%  "

# Row 8
$ws.Range("D8").Value = "'47.72"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.05%  "

# Row 9
$ws.Range("E9").Value = "  -0.35%  "

# Row 10
$ws.Range("D10").Value = "'1.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.17%  "

# Row 11
$ws.Range("D11").Value = "'0.07502"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.66%  "

# Row 13
$ws.Range("D13").Value = "'20.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.09%  "

# Row 14
$ws.Range("D14").Value = "'5.930"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.24%  "

# Row 15
$ws.Range("D15").Value = "'6.894"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.30%  "

# Row 16
$ws.Range("D16").Value = "1.558.64"
$ws.Range("E16").Value = "  -0.16%  "

# Row 17
$ws.Range("D17").Value = "'0.00001114"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.22%  "

# Row 18
$ws.Range("E18").Value = "  +0.26%  "

# Row 19
$ws.Range("D19").Value = "'0.06735"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.78%  "

# Row 20
$ws.Range("D20").Value = "'6.429"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.64%  "

# Row 21
$ws.Range("E21").Value = "  -0.07%  "

# Row 22
$ws.Range("D22").Value = "'16.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.60%  "

# Row 23
$ws.Range("D23").Value = "'11.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.04%  "

# Row 24
$ws.Range("D24").Value = "22.412.82"
$ws.Range("E24").Value = "  +0.03%  "

# Row 25
$ws.Range("E25").Value = "  -1.82%  "

# Row 26
$ws.Range("D26").Value = "'2.627"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.24%  "

# Row 27
$ws.Range("D27").Value = "'150.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.63%  "

# Row 28
$ws.Range("D28").Value = "'19.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.10%  "

# Row 29
$ws.Range("D29").Value = "'4.942"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.20%  "

# Row 30
$ws.Range("D30").Value = "'124.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.60%  "

# Row 31
$ws.Range("D31").Value = "1.736.61"
$ws.Range("E31").Value = "  +0.07%  "

# Row 32
$ws.Range("D32").Value = "'1.079"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.33%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'6.079"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.45%  "

# Row 34
$ws.Range("B34").Value = "WEMIXTOKEN"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'1.979"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.89%  "

# Row 35
$ws.Range("D35").Value = "'9.865"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.99%  "

# Row 36
$ws.Range("E36").Value = "  +0.43%  "

# Row 37
$ws.Range("E37").Value = "  +2.01%  "

# Row 38
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.06378"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.04%  "

# Row 39
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.298"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.20%  "

# Row 40
$ws.Range("D40").Value = "'0.2207"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.10%  "

# Row 41
$ws.Range("D41").Value = "'5.320"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.33%  "

# Row 42
$ws.Range("D42").Value = "'11.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.57%  "

# Row 43
$ws.Range("D43").Value = "'0.6229"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.00%  "

# Row 44
$ws.Range("E44").Value = "  +0.02%  "

# Row 45
$ws.Range("D45").Value = "'13.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.85%  "

# Row 46
$ws.Range("D46").Value = "'0.6038"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.18%  "

# Row 47
$ws.Range("D47").Value = "'3.774"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.70%  "

# Row 48
$ws.Range("E48").Value = "  +1.70%  "

# Row 49
$ws.Range("D49").Value = "'124.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.32%  "

# Row 50
$ws.Range("D50").Value = "'1.192"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.48%  "

# Row 51
$ws.Range("D51").Value = "'0.07189"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.16%  "
